# Weekly update: insert the newest weekly price record for
# "Vega Monumental Concepción - Mango" at the top of the data block
# (row 143), pushing the existing rows 143:155 down to 144:156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 143; rows 143:155 shift down to 144:156.
$ws.Rows.Item(143).Insert()

# Populate the new row 143 with this week's record.
$ws.Range("A143").Value = 11
$ws.Range("B143").Value = "Vega Monumental Concepción"
$ws.Range("C143").Value = "Bíobío"
$ws.Range("D143").Value = 45013
$ws.Range("E143").Value = 8
$ws.Range("F143").Value = "Fruta"
$ws.Range("G143").Value = 100108
$ws.Range("H143").Value = "Tropicales y subtropicales"
$ws.Range("I143").Value = 100108002
$ws.Range("J143").Value = "Mango"
$ws.Range("K143").Value = "Sin especificar"
$ws.Range("L143").Value = "Primera"
$ws.Range("M143").Value = 200
$ws.Range("N143").Value = 7500
$ws.Range("O143").Value = 8000
$ws.Range("P143").Value = 7750
$ws.Range("Q143").Value = "$/bandeja 4 kilos"
$ws.Range("R143").Value = "Perú"
$ws.Range("S143").Value = 1938
$ws.Range("T143").Value = 4
